$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = "aa"
$ws.Range("J14").Value = "Agree/Accept"
$ws.Range("I19").Value = "%"
$ws.Range("J19").Value = "Uninterpretable"
$ws.Range("I26").Value = "%"
$ws.Range("J26").Value = "Uninterpretable"
$ws.Range("I53").Value = "aa"
$ws.Range("J53").Value = "Agree/Accept"
$ws.Range("I66").Value = "aa"
$ws.Range("J66").Value = "Agree/Accept"
$ws.Range("I67").Value = "b"
$ws.Range("J67").Value = "Acknowledge (Backchannel)"
$ws.Range("I73").Value = "aa"
$ws.Range("J73").Value = "Agree/Accept"
$ws.Range("I76").Value = "sd"
$ws.Range("J76").Value = "Statement-non-opinion"
$ws.Range("I85").Value = "aa"
$ws.Range("J85").Value = "Agree/Accept"
$ws.Range("I90").Value = "sd"
$ws.Range("J90").Value = "Statement-non-opinion"
$ws.Range("I104").Value = "sd"
$ws.Range("J104").Value = "Statement-non-opinion"
$ws.Range("I112").Value = "aa"
$ws.Range("J112").Value = "Agree/Accept"
$ws.Range("I114").Value = "sd"
$ws.Range("J114").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sv"
$ws.Range("J132").Value = "Statement-opinion"
$ws.Range("I137").Value = "b"
$ws.Range("J137").Value = "Acknowledge (Backchannel)"
$ws.Range("I138").Value = "ba"
$ws.Range("J138").Value = "Appreciation"
$ws.Range("I142").Value = "sd"
$ws.Range("J142").Value = "Statement-non-opinion"
$ws.Range("I145").Value = "aa"
$ws.Range("J145").Value = "Agree/Accept"
$ws.Range("I157").Value = "%"
$ws.Range("J157").Value = "Uninterpretable"
$ws.Range("I160").Value = "sv"
$ws.Range("J160").Value = "Statement-opinion"
$ws.Range("I161").Value = "aa"
$ws.Range("J161").Value = "Agree/Accept"
$ws.Range("I163").Value = "sv"
$ws.Range("J163").Value = "Statement-opinion"
$ws.Range("I167").Value = "b"
$ws.Range("J167").Value = "Acknowledge (Backchannel)"
$ws.Range("I177").Value = "sd"
$ws.Range("J177").Value = "Statement-non-opinion"
$ws.Range("I178").Value = "b"
$ws.Range("J178").Value = "Acknowledge (Backchannel)"
$ws.Range("I188").Value = "ba"
$ws.Range("J188").Value = "Appreciation"
$ws.Range("I190").Value = "b"
$ws.Range("J190").Value = "Acknowledge (Backchannel)"
$ws.Range("I193").Value = "sd"
$ws.Range("J193").Value = "Statement-non-opinion"
$ws.Range("I195").Value = "sv"
$ws.Range("J195").Value = "Statement-opinion"
$ws.Range("I200").Value = "sv"
$ws.Range("J200").Value = "Statement-opinion"
$ws.Range("I205").Value = "sv"
$ws.Range("J205").Value = "Statement-opinion"
$ws.Range("I209").Value = "sv"
$ws.Range("J209").Value = "Statement-opinion"
$ws.Range("I220").Value = "sd"
$ws.Range("J220").Value = "Statement-non-opinion"
$ws.Range("I227").Value = "sv"
$ws.Range("J227").Value = "Statement-opinion"
$ws.Range("I228").Value = "%"
$ws.Range("J228").Value = "Uninterpretable"
$ws.Range("I234").Value = "sv"
$ws.Range("J234").Value = "Statement-opinion"
$ws.Range("I236").Value = "sd"
$ws.Range("J236").Value = "Statement-non-opinion"
$ws.Range("I239").Value = "sd"
$ws.Range("J239").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "b"
$ws.Range("J243").Value = "Acknowledge (Backchannel)"
